$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the last row (old row 22, "Bibliografia:" / long bibliography text).
#    This also shrinks the used range / dimension automatically.
$ws.Rows.Item(22).Delete()

# 2) Rewrite the label column (A) and the two content columns (B/C) for rows
#    10-21 to match the new layout (rows shifted/re-paired around the
#    "Docentes responsaveis" / teacher-name block being relocated).

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "8711686 - Flavia Reis Cardoso Rojas"
$ws.Range("C10").Value = "8711686 - Flavia Reis Cardoso Rojas"

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Presenting to the students the introductory concepts of Physics and in particular, of Mechanics including kinematics and dynamics, including basic concepts of statistical and data analysis."
$ws.Range("C11").Value = "Presenting to the students the introductory concepts of Physics and in particular, of Mechanics including kinematics and dynamics, including basic concepts of statistical and data analysis."

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular."
$ws.Range("C14").Value = "Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular."

$ws.Range("A15").Value = "Programa:"
# "01/01/2018" would otherwise be auto-recognised as a date; force it back
# to plain text (like the rest of the sheet) and restore the shared column
# formatting afterwards.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2018"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton's laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum."
$ws.Range("C16").Value = "1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton's laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum."

$ws.Range("A17").Value = "Avaliação:"
# Row 17 only carries the label now - drop its old B/C content entirely.
$ws.Range("B17:C17").Clear()

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8711686 - Flavia Reis Cardoso Rojas"
$ws.Range("C18").Value = "8711686 - Flavia Reis Cardoso Rojas"
# B18/C18 are brand-new cells in this row; pick up the same number format /
# font / alignment used by the rest of column B and C instead of leaving
# them on the default style.
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada"
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada"

# 3) Row heights for the re-shuffled rows.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(9).RowHeight
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
